$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - F column updates ("想去人数" / want-to-go count)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 8292
$wsExhibit.Range("F3").Value = 7700
$wsExhibit.Range("F4").Value = 114
$wsExhibit.Range("F9").Value = 110
$wsExhibit.Range("F10").Value = 156
$wsExhibit.Range("F11").Value = 226
$wsExhibit.Range("F12").Value = 697
$wsExhibit.Range("F13").Value = 123
$wsExhibit.Range("F14").Value = 1262
$wsExhibit.Range("F16").Value = 47
$wsExhibit.Range("F17").Value = 9
$wsExhibit.Range("F19").Value = 106

# Sheet "全部类型" (All Types) - F column updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 8292
$wsAll.Range("F3").Value = 7700
$wsAll.Range("F4").Value = 114
$wsAll.Range("F9").Value = 111
$wsAll.Range("F10").Value = 156
$wsAll.Range("F11").Value = 226
$wsAll.Range("F12").Value = 697
$wsAll.Range("F13").Value = 123
$wsAll.Range("F14").Value = 1262
$wsAll.Range("F16").Value = 47
$wsAll.Range("F17").Value = 9
$wsAll.Range("F19").Value = 106
